# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (shape "Google Shape;122;p17") gets a new built-in
#    table style applied (tableStyleId GUID change).
# 2) The deck's theme colour palette is changed from the "Integral" /
#    "Red Violet" palette to the stock "Office Theme" / "Office" palette
#    (the RGB values that back the slide master's theme part).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 --------------------------------
$slide  = $p.Slides.Item(5)
$shape  = $slide.Shapes.Item(2)
$table  = $shape.Table
$table.ApplyStyle("{8751A8B8-C712-4FC7-916E-42FA18597F95}")

# --- 2. Swap the active theme's colour scheme over to "Office" -------
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Item(1).RGB  = 0          # dk1      000000
$themeColors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388    # dk2      44546A
$themeColors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407      # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Item(10).RGB = 4697456    # accent6  70AD47
$themeColors.Item(11).RGB = 12673797   # hlink    0563C1
$themeColors.Item(12).RGB = 7491477    # folHlink 954F72
